$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new exposure sites were added at the top of the list (rows 2 and 3),
# pushing every existing data row down by two positions (old row N -> new
# row N+2). Shift the existing data down manually, working from the bottom
# of the sheet upward so we never overwrite a row before it has been read.
for ($r = 160; $r -ge 2; $r--) {
    $dest = $r + 2
    $ws.Cells.Item($dest, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($dest, 2).Value = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($dest, 3).Value = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($dest, 4).Value = $ws.Cells.Item($r, 4).Value2
}

# Populate the two newly freed-up rows with the new exposure site data.
$ws.Range("A2").Value = "185 Cooper St, Epping VIC 3076"
$ws.Range("B2").Value = -37.653023
$ws.Range("C2").Value = 145.014685
$ws.Range("D2").Value = "Whittlesea (C)"

$ws.Range("A3").Value = "260-264 Arthur St, Fairfield VIC 3078"
$ws.Range("B3").Value = -37.766715
$ws.Range("C3").Value = 145.021149
$ws.Range("D3").Value = "Darebin (C)"
